$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.210.57'
$ws.Range("E2").Value = '  +1.80%  '

# Row 3
$ws.Range("D3").Value = '2.570.05'
$ws.Range("E3").Value = '  +1.81%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.37'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.18%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.05'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +2.55%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.575'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +0.37%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -0.13%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +3.64%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.82'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +0.09%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.88%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.50'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -0.31%  '

# Row 13
$ws.Range("E13").Value = '  -4.64%  '

# Row 14
$ws.Range("D14").Value = '2.958.45'
$ws.Range("E14").Value = '  +1.57%  '

# Row 15
$ws.Range("D15").Value = '2.528.85'
$ws.Range("E15").Value = '  +0.96%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.11'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -1.83%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +0.63%  '

# Row 18
$ws.Range("D18").Value = '43.142.01'
$ws.Range("E18").Value = '  +1.57%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.89'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +5.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.66'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -2.57%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0965'
$ws.Range("E21").Value = '  +0.68%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.09'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +0.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.73'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +1.74%  '

# Row 24
$ws.Range("E24").Value = '  -0.07%  '

# Row 25
$ws.Range("E25").Value = '  +2.49%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.83'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +1.73%  '

# Row 27
$ws.Range("E27").Value = '  +0.42%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +2.10%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.42'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +3.89%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.31'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +1.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.84'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -1.64%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.97'
$ws.Range("D32").NumberFormat = "General"

# Row 33
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.37'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +3.13%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.14'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +1.56%  '

# Row 35
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.28'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +0.70%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0808'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +3.27%  '

# Row 37
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.70'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +2.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.113'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +2.10%  '

# Row 39
$ws.Range("E39").Value = '  +4.46%  '

# Row 40
$ws.Range("E40").Value = '  +0.22%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.24'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -6.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.88'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +2.62%  '

# Row 43
$ws.Range("E43").Value = '  +1.78%  '

# Row 44
$ws.Range("E44").Value = '  -0.02%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.26'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -0.05%  '

# Row 46
$ws.Range("D46").Value = '1.987.65'
$ws.Range("E46").Value = '  -0.87%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.88'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +0.18%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.05'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +2.83%  '

# Row 49
$ws.Range("D49").Value = '2.806.32'
$ws.Range("E49").Value = '  +1.63%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.66'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +2.70%  '

# Row 51
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.40'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +0.05%  '
